$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "NA" survival/size values for columns W (day19) and X (day20)
# across rows 29-36, completing the final set of survival and size data.
$ws.Range("W29:X36").Value = "NA"

# Update the sheet view to reflect the newly active/selected range
$excel.ActiveWindow.ScrollColumn = 18
$ws.Range("W29:X36").Select()
